$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows for d=1..d=5, d=7 and d=10 (rows 2-8).
# We need to insert a new "d=6" row between the existing "d=5" row (row 6)
# and "d=7" row (row 7), pushing "d=7" and "d=10" down by one row.
#
# Shift the existing rows 7 ("d=7") and 8 ("d=10") down to rows 8 and 9
# using Copy(Destination), which carries over both values and formatting
# (so the bold/centered/bordered label style is preserved) without
# introducing new, unused cell styles like Rows.Insert() would.

# Move "d=10" row (row 8) down to row 9.
$ws.Range("A8:E8").Copy($ws.Range("A9:E9"))

# Move "d=7" row (row 7) down to row 8.
$ws.Range("A7:E7").Copy($ws.Range("A8:E8"))

# Write the new "d=6" row into row 7 with its corrected results.
$ws.Range("A7").Value = "d=6"
$ws.Range("B7").Value = 97.68344439838881
$ws.Range("C7").Value = 97.75535123068441
$ws.Range("D7").Value = 97.79188381272378
$ws.Range("E7").Value = 97.75022338235779

$wb.Save()
